$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 / D3: the long "Tabla perforata ... De la 1 la 4" note moves from row 2
# down to row 3 (and loses its trailing space).
$ws.Range("D2").Value = "Tabla perforata vopsita electrostatic"
$ws.Range("D3").Value = "Tabla perforata vopsita electrostatic. De la 1 la 4"

# B11: the "Poza" reference for this row now points at image7 instead of
# the (now unused) image9.
$ws.Range("B11").Value = "assets/image7.jpg"

# C12 / C14: the area callouts gain a "suprafata peretilor" prefix and the
# cells switch to wrapped text.
$ws.Range("C12").Value = "suprafata peretilor     170 m²"
$ws.Range("C12").WrapText = $true

$ws.Range("C14").Value = "suprafata peretilor        178.9 m²"
$ws.Range("C14").WrapText = $true

# Restore the scroll position / active selection recorded for this sheet.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I12").Select()
